$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Headers for new columns I and J ---
# Copy the formatting of the existing header style (H1) onto the new header
# cells before putting in their text, so they share the same style index
# (bold, centered, bordered) as the other header cells.
$ws.Range("H1").Copy()
$ws.Range("I1:J1").PasteSpecial(-4122)
$excel.CutCopyMode = 0

$ws.Range("I1").Value = "I0"
$ws.Range("J1").Value = "IF"

# --- Data values for I2:I77 and J2:J77 ---
$iValues = @(7,7,7,5,9,7,8,6,7,9,3,5,8,3,8,4,11,6,7,5,7,6,7,8,7,7,8,9,5,6,7,8,6,6,7,7,7,8,8,4,6,3,6,7,9,9,7,7,6,8,10,5,7,6,9,6,7,5,4,8,8,7,9,9,8,9,9,5,7,9,5,1,7,5,6,4)
$jValues = @(7,7,7,6,9,7,8,6,7,9,3,6,8,4,8,5,11,7,7,6,8,7,7,8,7,7,8,9,5,7,7,8,6,6,7,7,7,8,8,4,6,4,7,7,9,9,8,8,6,8,10,6,7,6,9,7,7,6,5,8,8,8,9,9,8,9,9,6,7,9,5,1,7,5,6,4)

for ($r = 2; $r -le 77; $r++) {
    $idx = $r - 2
    $ws.Cells.Item($r, 9).Value = $iValues[$idx]
    $ws.Cells.Item($r, 10).Value = $jValues[$idx]
}
